# Insert a new data row right before the current row 311 (Feria Lagunitas
# de Puerto Montt - Coliflor), shifting the existing rows 311-327 down to
# 312-328, and fill the new row with its own values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(311).Insert()

$ws.Cells.Item(311, 1).Value = 4
$ws.Cells.Item(311, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(311, 3).Value = "Los Lagos"
$ws.Cells.Item(311, 4).Value = 44706
$ws.Cells.Item(311, 5).Value = 10
$ws.Cells.Item(311, 6).Value = 100112008
$ws.Cells.Item(311, 7).Value = "Coliflor"
$ws.Cells.Item(311, 8).Value = "Sin especificar"
$ws.Cells.Item(311, 9).Value = "Primera"
$ws.Cells.Item(311, 10).Value = 250
$ws.Cells.Item(311, 11).Value = 1800
$ws.Cells.Item(311, 12).Value = 1800
$ws.Cells.Item(311, 13).Value = 1800
$ws.Cells.Item(311, 14).Value = "`$/unidad"
$ws.Cells.Item(311, 15).Value = "Región del Maule"
$ws.Cells.Item(311, 16).Value = 1800
$ws.Cells.Item(311, 17).Value = 1
$ws.Cells.Item(311, 18).Value = "Hortaliza"
